# The workbook tracks daily Arándano (blue) price observations at the
# Lo Valledor wholesale market. Two new daily records need to be inserted
# at the top of the data block (rows 639-640), pushing the existing
# records (formerly rows 639-715) down by two rows to 641-717.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block; this shifts rows
# 639:715 down to 641:717 and keeps row/column formatting consistent.
$ws.Rows("639:640").Insert()

# --- New row 639 ---
$ws.Cells.Item(639, 1).Value = 6
$ws.Cells.Item(639, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(639, 3).Value = "Metropolitana"
$ws.Cells.Item(639, 4).Value = 45154
$ws.Cells.Item(639, 5).Value = 13
$ws.Cells.Item(639, 6).Value = "Fruta"
$ws.Cells.Item(639, 7).Value = 100101
$ws.Cells.Item(639, 8).Value = "Berries"
$ws.Cells.Item(639, 9).Value = 100101001
$ws.Cells.Item(639, 10).Value = "Arándano (blue)"
$ws.Cells.Item(639, 11).Value = "Sin especificar"
$ws.Cells.Item(639, 12).Value = "Especial"
$ws.Cells.Item(639, 13).Value = 250
$ws.Cells.Item(639, 14).Value = 18000
$ws.Cells.Item(639, 15).Value = 18000
$ws.Cells.Item(639, 16).Value = 18000
$ws.Cells.Item(639, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(639, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(639, 19).Value = 9000
$ws.Cells.Item(639, 20).Value = 2

# --- New row 640 ---
$ws.Cells.Item(640, 1).Value = 6
$ws.Cells.Item(640, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(640, 3).Value = "Metropolitana"
$ws.Cells.Item(640, 4).Value = 45154
$ws.Cells.Item(640, 5).Value = 13
$ws.Cells.Item(640, 6).Value = "Fruta"
$ws.Cells.Item(640, 7).Value = 100101
$ws.Cells.Item(640, 8).Value = "Berries"
$ws.Cells.Item(640, 9).Value = 100101001
$ws.Cells.Item(640, 10).Value = "Arándano (blue)"
$ws.Cells.Item(640, 11).Value = "Sin especificar"
$ws.Cells.Item(640, 12).Value = "Primera"
$ws.Cells.Item(640, 13).Value = 1070
$ws.Cells.Item(640, 14).Value = 12000
$ws.Cells.Item(640, 15).Value = 12000
$ws.Cells.Item(640, 16).Value = 12000
$ws.Cells.Item(640, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(640, 18).Value = "Perú"
$ws.Cells.Item(640, 19).Value = 8000
$ws.Cells.Item(640, 20).Value = 1.5
